$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 304.2857
$ws.Range("J9").Value = 600
$ws.Range("L9").Value = 600
$ws.Range("N9").Value = -938

$ws.Range("H17").Value = 3230531.8
$ws.Range("J17").Value = 3576543.8
$ws.Range("L17").Value = 10729631.4
$ws.Range("N17").Value = -10729967.4

$ws.Range("H51").Value = 2496.3635
$ws.Range("I51").Value = 2625
$ws.Range("J51").Value = 2422.8572
$ws.Range("K51").Value = 2625
$ws.Range("L51").Value = 2422.8572
$ws.Range("M51").Value = -2141
$ws.Range("N51").Value = -3390.8572

$ws.Range("H98").Value = 927
$ws.Range("I98").Value = 927
$ws.Range("K98").Value = 927
$ws.Range("M98").Value = 571

$ws.Range("H116").Value = 5020.5
$ws.Range("I116").Value = 1200
$ws.Range("J116").Value = 5445
$ws.Range("K116").Value = 1200
$ws.Range("L116").Value = 5445
$ws.Range("N116").Value = -12329
$ws.Range("M116").Value = 2242

$ws.Range("H122").Value = 927
$ws.Range("I122").Value = 927
$ws.Range("K122").Value = 2781
$ws.Range("M122").Value = -331

$ws.Range("H129").Value = 165033.14
$ws.Range("J129").Value = 182982.36
$ws.Range("L129").Value = 548947.08
$ws.Range("N129").Value = -558947.08

$ws.Range("H137").Value = 1690.8182
$ws.Range("I137").Value = 1562.375
$ws.Range("J137").Value = 2033.3334
$ws.Range("K137").Value = 4687.125
$ws.Range("L137").Value = 6100.0002
$ws.Range("M137").Value = -2137.125
$ws.Range("N137").Value = -11200.0002

$ws.Range("H138").Value = 1613.186
$ws.Range("I138").Value = 594.20514
$ws.Range("J138").Value = 2458.7234
$ws.Range("K138").Value = 1782.61542
$ws.Range("L138").Value = 7376.1702
$ws.Range("M138").Value = 3357.38458
$ws.Range("N138").Value = -17656.1702

$ws.Range("H141").Value = 1278.4615
$ws.Range("I141").Value = 1169.6
$ws.Range("K141").Value = 3508.8
$ws.Range("M141").Value = 1671.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4811.837
$ws.Range("I32").Value = 5172.4653
$ws.Range("J32").Value = 2227.3333
$ws.Range("K32").Value = 5172.4653
$ws.Range("L32").Value = 2227.3333
$ws.Range("M32").Value = -4885.4653
$ws.Range("N32").Value = -2801.3333

$ws.Range("H45").Value = 4045.5715
$ws.Range("J45").Value = 4169
$ws.Range("L45").Value = 4169
$ws.Range("N45").Value = -4923

$ws.Range("H74").Value = 37039692
$ws.Range("I74").Value = 50002660
$ws.Range("K74").Value = 50002660
$ws.Range("M74").Value = -50001786

$ws.Range("H77").Value = 37039692
$ws.Range("I77").Value = 50002660
$ws.Range("K77").Value = 250013300
$ws.Range("M77").Value = -250008932

$ws.Range("H102").Value = 1344.2858
$ws.Range("I102").Value = 803.3333
$ws.Range("K102").Value = 803.3333
$ws.Range("M102").Value = 818.6667

$ws.Range("H122").Value = 5198.625
$ws.Range("I122").Value = 4531.6665
$ws.Range("J122").Value = 7199.5
$ws.Range("K122").Value = 13594.9995
$ws.Range("L122").Value = 21598.5
$ws.Range("M122").Value = -11144.9995
$ws.Range("N122").Value = -26498.5

$ws.Range("H132").Value = 15358.243
$ws.Range("I132").Value = 1915.8518
$ws.Range("K132").Value = 5747.555399999999
$ws.Range("M132").Value = -3217.555399999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3280.0977
$ws.Range("I134").Value = 3555.4412
$ws.Range("J134").Value = 1942.7142
$ws.Range("K134").Value = 10666.3236
$ws.Range("L134").Value = 5828.142599999999
$ws.Range("M134").Value = -8131.3236
$ws.Range("N134").Value = -10898.1426

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 19275.105
$ws.Range("I31").Value = 44030.43
$ws.Range("K31").Value = 44030.43
$ws.Range("M31").Value = -43735.43

$ws.Range("H34").Value = 19275.105
$ws.Range("I34").Value = 44030.43
$ws.Range("K34").Value = 44030.43
$ws.Range("M34").Value = -43828.43

$ws.Range("H58").Value = 24060.092
$ws.Range("I58").Value = 1391.6428
$ws.Range("J58").Value = 63729.875
$ws.Range("K58").Value = 1391.6428
$ws.Range("L58").Value = 63729.875
$ws.Range("M58").Value = -1188.6428
$ws.Range("N58").Value = -64135.875

$ws.Range("H105").Value = 31251732
$ws.Range("I105").Value = 41667304
$ws.Range("K105").Value = 41667304
$ws.Range("M105").Value = -41665557

$ws.Range("H122").Value = 1682.3043
$ws.Range("I122").Value = 2025.3636
$ws.Range("K122").Value = 6076.0908
$ws.Range("M122").Value = -3626.0908

$ws.Range("H132").Value = 20545.172
$ws.Range("I132").Value = 26321.762
$ws.Range("J132").Value = 5381.625
$ws.Range("K132").Value = 78965.28599999999
$ws.Range("L132").Value = 16144.875
$ws.Range("M132").Value = -76435.28599999999
$ws.Range("N132").Value = -21204.875

$ws.Range("H136").Value = 24060.092
$ws.Range("I136").Value = 1391.6428
$ws.Range("J136").Value = 63729.875
$ws.Range("K136").Value = 4174.928400000001
$ws.Range("L136").Value = 191189.625
$ws.Range("M136").Value = -1624.928400000001
$ws.Range("N136").Value = -196289.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 2875.5386
$ws.Range("I3").Value = 1331.7273
$ws.Range("J3").Value = 11366.5
$ws.Range("K3").Value = 3995.1819
$ws.Range("L3").Value = 34099.5
$ws.Range("M3").Value = -3883.1819
$ws.Range("N3").Value = -34323.5

$ws.Range("H122").Value = 515.7143
$ws.Range("I122").Value = 250.23077
$ws.Range("J122").Value = 947.125
$ws.Range("K122").Value = 2252.07693
$ws.Range("L122").Value = 8524.125
$ws.Range("M122").Value = 197.9230699999998
$ws.Range("N122").Value = -13424.125

$ws.Range("H129").Value = 278685.94
$ws.Range("I129").Value = 776.6667
$ws.Range("J129").Value = 334267.8
$ws.Range("K129").Value = 2330.0001
$ws.Range("L129").Value = 1002803.4
$ws.Range("M129").Value = 2669.9999
$ws.Range("N129").Value = -1012803.4

$ws.Range("H131").Value = 806.34
$ws.Range("I131").Value = 703
$ws.Range("J131").Value = 809.5361
$ws.Range("K131").Value = 2109
$ws.Range("L131").Value = 2428.6083
$ws.Range("M131").Value = 2931
$ws.Range("N131").Value = -12508.6083

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H92").Value = 8625
$ws.Range("J92").Value = 8625
$ws.Range("L92").Value = 8625
$ws.Range("N92").Value = -12369

$ws.Range("H93").Value = 13999.833
$ws.Range("J93").Value = 13999.833
$ws.Range("L93").Value = 13999.833
$ws.Range("N93").Value = -17743.833

$ws.Range("H102").Value = 26317044
$ws.Range("I102").Value = 31251204
$ws.Range("J102").Value = 1526.6666
$ws.Range("K102").Value = 31251204
$ws.Range("L102").Value = 1526.6666
$ws.Range("M102").Value = -31249582
$ws.Range("N102").Value = -4770.6666

$ws.Range("H122").Value = 51283910
$ws.Range("I122").Value = 23810682
$ws.Range("J122").Value = 83336020
$ws.Range("K122").Value = 71432046
$ws.Range("L122").Value = 250008060
$ws.Range("M122").Value = -71429596
$ws.Range("N122").Value = -250012960

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3282.6667
$ws.Range("I22").Value = 3099.111
$ws.Range("J22").Value = 3833.3333
$ws.Range("K22").Value = 3099.111
$ws.Range("L22").Value = 3833.3333
$ws.Range("M22").Value = -2804.111
$ws.Range("N22").Value = -4423.3333

$ws.Range("H27").Value = 3282.6667
$ws.Range("I27").Value = 3099.111
$ws.Range("J27").Value = 3833.3333
$ws.Range("K27").Value = 3099.111
$ws.Range("L27").Value = 3833.3333
$ws.Range("M27").Value = -2992.111
$ws.Range("N27").Value = -4047.3333

$ws.Range("H40").Value = 7084.3335
$ws.Range("I40").Value = 5199.4
$ws.Range("K40").Value = 5199.4
$ws.Range("M40").Value = -5063.4

$ws.Range("H111").Value = 33591.332
$ws.Range("J111").Value = 33591.332
$ws.Range("L111").Value = 33591.332
$ws.Range("N111").Value = -41771.332

$ws.Range("H122").Value = 855265
$ws.Range("J122").Value = 2666.4375
$ws.Range("L122").Value = 7999.3125
$ws.Range("N122").Value = -12899.3125

$ws.Range("H132").Value = 2242.25
$ws.Range("J132").Value = 4166.3335
$ws.Range("L132").Value = 12499.0005
$ws.Range("N132").Value = -17559.0005

$ws.Range("H136").Value = 4000
$ws.Range("I136").Value = 4000
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 12000
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -9450
$ws.Range("N136").Value = $null

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1501.9565
$ws.Range("I122").Value = 1573.4736
$ws.Range("J122").Value = 1162.25
$ws.Range("K122").Value = 4720.4208
$ws.Range("L122").Value = 3486.75
$ws.Range("M122").Value = -2270.4208
$ws.Range("N122").Value = -8386.75

$ws.Range("H126").Value = 1428.5714
$ws.Range("I126").Value = 1375
$ws.Range("K126").Value = 4125
$ws.Range("M126").Value = -1655

$ws.Range("H132").Value = 3033
$ws.Range("I132").Value = 1500
$ws.Range("J132").Value = 3339.6
$ws.Range("K132").Value = 4500
$ws.Range("L132").Value = 10018.8
$ws.Range("M132").Value = -1970
$ws.Range("N132").Value = -15078.8

Write-Host "Applied all Typhon_Profits updates"
